$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 25329.5
$ws.Range("J3").Value = 25329.5
$ws.Range("L3").Value = 25329.5
$ws.Range("N3").Value = -25557.5
$ws.Range("H40").Value = 2408.818
$ws.Range("J40").Value = 2380.9375
$ws.Range("L40").Value = 2380.9375
$ws.Range("N40").Value = -2730.9375
$ws.Range("H69").Value = 16829.62
$ws.Range("J69").Value = 20178.215
$ws.Range("L69").Value = 60534.645
$ws.Range("N69").Value = -62282.645
$ws.Range("H72").Value = 16829.62
$ws.Range("J72").Value = 20178.215
$ws.Range("L72").Value = 181603.935
$ws.Range("N72").Value = -190339.935
$ws.Range("H94").Value = 441.3
$ws.Range("I94").Value = 251.625
$ws.Range("K94").Value = 251.625
$ws.Range("M94").Value = 199.375
$ws.Range("H102").Value = 25329.5
$ws.Range("J102").Value = 25329.5
$ws.Range("L102").Value = 25329.5
$ws.Range("N102").Value = -31819.5
$ws.Range("H106").Value = 146668670
$ws.Range("I106").Value = 220000500
$ws.Range("J106").Value = 5000
$ws.Range("K106").Value = 220000500
$ws.Range("L106").Value = 5000
$ws.Range("M106").Value = -219999869
$ws.Range("N106").Value = -6262
$ws.Range("H135").Value = 987.67645
$ws.Range("I135").Value = 821.1724
$ws.Range("K135").Value = 7390.551600000001
$ws.Range("M135").Value = -4855.551600000001
$ws.Range("H137").Value = 1956.7709
$ws.Range("I137").Value = 1734.2821
$ws.Range("K137").Value = 5202.846299999999
$ws.Range("M137").Value = -2652.846299999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5800
$ws.Range("I61").Value = 5720.841
$ws.Range("K61").Value = 5720.841
$ws.Range("M61").Value = -5508.841
$ws.Range("H74").Value = 9436363
$ws.Range("I74").Value = 12196645
$ws.Range("J74").Value = 5398.75
$ws.Range("K74").Value = 12196645
$ws.Range("L74").Value = 5398.75
$ws.Range("M74").Value = -12195771
$ws.Range("N74").Value = -7146.75
$ws.Range("H77").Value = 9436363
$ws.Range("I77").Value = 12196645
$ws.Range("J77").Value = 5398.75
$ws.Range("K77").Value = 60983225
$ws.Range("L77").Value = 26993.75
$ws.Range("M77").Value = -60978857
$ws.Range("N77").Value = -35729.75
$ws.Range("H102").Value = 2876.1875
$ws.Range("I102").Value = 2835
$ws.Range("K102").Value = 2835
$ws.Range("M102").Value = -1213
$ws.Range("H134").Value = 82237.8
$ws.Range("J134").Value = 92797.25
$ws.Range("L134").Value = 92797.25
$ws.Range("N134").Value = -102937.25
$ws.Range("H136").Value = 5800
$ws.Range("I136").Value = 5720.841
$ws.Range("K136").Value = 17162.523
$ws.Range("M136").Value = -14612.523

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 66427.71000000001
$ws.Range("I59").Value = 30000
$ws.Range("K59").Value = 30000
$ws.Range("M59").Value = -29153
$ws.Range("H86").Value = 3849.7827
$ws.Range("I86").Value = 2971.75
$ws.Range("K86").Value = 2971.75
$ws.Range("M86").Value = -1848.75
$ws.Range("H89").Value = 3849.7827
$ws.Range("I89").Value = 2971.75
$ws.Range("K89").Value = 14858.75
$ws.Range("M89").Value = -9242.75
$ws.Range("H105").Value = 853833.9
$ws.Range("I105").Value = 1083543.4
$ws.Range("K105").Value = 1083543.4
$ws.Range("M105").Value = -1081796.4
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H134").Value = 2787.5
$ws.Range("I134").Value = 2089.2
$ws.Range("K134").Value = 6267.599999999999
$ws.Range("M134").Value = -3732.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8133.2
$ws.Range("I31").Value = 2579.6924
$ws.Range("J31").Value = 14149.5
$ws.Range("K31").Value = 2579.6924
$ws.Range("L31").Value = 14149.5
$ws.Range("M31").Value = -2284.6924
$ws.Range("N31").Value = -14739.5
$ws.Range("H34").Value = 8133.2
$ws.Range("I34").Value = 2579.6924
$ws.Range("J34").Value = 14149.5
$ws.Range("K34").Value = 2579.6924
$ws.Range("L34").Value = 14149.5
$ws.Range("M34").Value = -2377.6924
$ws.Range("N34").Value = -14553.5
$ws.Range("H52").Value = 64344.855
$ws.Range("I52").Value = 71950
$ws.Range("K52").Value = 71950
$ws.Range("M52").Value = -71656
$ws.Range("H58").Value = 3696.276
$ws.Range("I58").Value = 2709.6843
$ws.Range("K58").Value = 2709.6843
$ws.Range("M58").Value = -2506.6843
$ws.Range("H122").Value = 3623
$ws.Range("I122").Value = 1996
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 5988
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -3538
$ws.Range("N122").Value = -20650
$ws.Range("H132").Value = 1615.2142
$ws.Range("I132").Value = 1318.1
$ws.Range("J132").Value = 4091.1667
$ws.Range("K132").Value = 3954.3
$ws.Range("L132").Value = 12273.5001
$ws.Range("M132").Value = -1424.3
$ws.Range("N132").Value = -17333.5001
$ws.Range("H134").Value = 1822.5758
$ws.Range("I134").Value = 1298.362
$ws.Range("J134").Value = 5623.125
$ws.Range("K134").Value = 3895.086
$ws.Range("L134").Value = 16869.375
$ws.Range("M134").Value = -1360.086
$ws.Range("N134").Value = -21939.375
$ws.Range("H135").Value = 103721.164
$ws.Range("J135").Value = 103721.164
$ws.Range("L135").Value = 103721.164
$ws.Range("N135").Value = -113861.164
$ws.Range("H136").Value = 3696.276
$ws.Range("I136").Value = 2709.6843
$ws.Range("K136").Value = 8129.0529
$ws.Range("M136").Value = -5579.0529
$ws.Range("H140").Value = 115866.47
$ws.Range("J140").Value = 115642.64
$ws.Range("L140").Value = 115642.64
$ws.Range("N140").Value = -126002.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3345307.2
$ws.Range("I113").Value = 547.8
$ws.Range("J113").Value = 5435782
$ws.Range("K113").Value = 1643.4
$ws.Range("L113").Value = 16307346
$ws.Range("M113").Value = 526.6000000000001
$ws.Range("N113").Value = -16311686

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4027.0322
$ws.Range("I126").Value = 3493.8262
$ws.Range("K126").Value = 10481.4786
$ws.Range("M126").Value = -8011.4786
$ws.Range("H132").Value = 2851.8809
$ws.Range("I132").Value = 2358.6287
$ws.Range("J132").Value = 5318.143
$ws.Range("K132").Value = 7075.886100000001
$ws.Range("L132").Value = 15954.429
$ws.Range("M132").Value = -4545.886100000001
$ws.Range("N132").Value = -21014.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4001.12
$ws.Range("I40").Value = 3524.2273
$ws.Range("K40").Value = 3524.2273
$ws.Range("M40").Value = -3388.2273
$ws.Range("H82").Value = 2141.4443
$ws.Range("I82").Value = 1320.25
$ws.Range("J82").Value = 2798.4
$ws.Range("K82").Value = 1320.25
$ws.Range("L82").Value = 2798.4
$ws.Range("M82").Value = -959.25
$ws.Range("N82").Value = -3520.4
$ws.Range("H85").Value = 2141.4443
$ws.Range("I85").Value = 1320.25
$ws.Range("J85").Value = 2798.4
$ws.Range("K85").Value = 1320.25
$ws.Range("L85").Value = 2798.4
$ws.Range("M85").Value = -72.25
$ws.Range("N85").Value = -5294.4
$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -51988
$ws.Range("H129").Value = 94500
$ws.Range("J129").Value = 94500
$ws.Range("L129").Value = 94500
$ws.Range("N129").Value = -104500
$ws.Range("H132").Value = 3546.7646
$ws.Range("I132").Value = 1781.8235
$ws.Range("J132").Value = 5311.706
$ws.Range("K132").Value = 5345.470499999999
$ws.Range("L132").Value = 15935.118
$ws.Range("M132").Value = -2815.470499999999
$ws.Range("N132").Value = -20995.118
$ws.Range("H140").Value = 75085.25
$ws.Range("J140").Value = 75085.25
$ws.Range("L140").Value = 75085.25
$ws.Range("N140").Value = -85445.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 94750
$ws.Range("I18").Value = 39500
$ws.Range("K18").Value = 39500
$ws.Range("M18").Value = -39327
$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -26988
$ws.Range("H138").Value = 49406.668
$ws.Range("J138").Value = 49406.668
$ws.Range("L138").Value = 49406.668
$ws.Range("N138").Value = -59686.668
